$p = $ppt.ActivePresentation

# 1. Update the fixed "datetimeFigureOut" date text from 11/1/2012 to
#    11/6/2012 on the slide master and every custom (slide) layout.
$m = $p.SlideMaster

for ($j = 1; $j -le $m.Shapes.Count; $j++) {
    $sh = $m.Shapes.Item($j)
    if ($sh.Name -like "Date Placeholder*") {
        $sh.TextFrame.TextRange.Text = "11/6/2012"
    }
}

for ($i = 1; $i -le $m.CustomLayouts.Count; $i++) {
    $l = $m.CustomLayouts.Item($i)
    for ($j = 1; $j -le $l.Shapes.Count; $j++) {
        $sh = $l.Shapes.Item($j)
        if ($sh.Name -like "Date Placeholder*") {
            $sh.TextFrame.TextRange.Text = "11/6/2012"
        }
    }
}

# 2. Remove the "http://vk.com/club33848893" link textbox from the
#    front (first) slide.
$s = $p.Slides.Item(1)
$sh = $s.Shapes.Item("TextBox 4")
$sh.Delete()
